# Update NATMI TPM output with new values (per commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (target cluster: ECs) - recomputed TPM values
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8174055
$ws.Range("N2").Value = 1.634811
$ws.Range("O2").Value = 0.09761389918309914
$ws.Range("P2").Value = 0.08793380245030721
$ws.Range("Q2").Value = 0.09909692968424999
$ws.Range("R2").Value = 0.396387718737
$ws.Range("S2").Value = 0.09761389918309914
$ws.Range("T2").Value = 0.08793380245030721

# Row 3 (target cluster: FAPs) - recomputed specificity values
$ws.Range("O3").Value = 0.2151446495736546
$ws.Range("P3").Value = 0.2907140367125996
$ws.Range("S3").Value = 0.2151446495736546
$ws.Range("T3").Value = 0.2907140367125996

# Row 4: target cluster changes from "MuSCs" to "Inflammatory-Mac", with new TPM values
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008574
$ws.Range("N4").Value = 0.025722
$ws.Range("O4").Value = 0.001023900098049123
$ws.Range("P4").Value = 0.001383544193565374
$ws.Range("Q4").Value = 0.001039456029
$ws.Range("R4").Value = 0.006236736173999999
$ws.Range("S4").Value = 0.001023900098049123
$ws.Range("T4").Value = 0.001383544193565374

# Row 5: target cluster changes from "Inflammatory-Mac" (formerly Neutrophils) to "MuSCs", with new TPM values
$ws.Range("D5").Value = "MuSCs"
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.712803
$ws.Range("N5").Value = 11.425606
$ws.Range("O5").Value = 0.6822182822294519
$ws.Range("P5").Value = 0.6145646077002448
$ws.Range("Q5").Value = 0.6925831025004999
$ws.Range("R5").Value = 2.770332410002
$ws.Range("S5").Value = 0.6822182822294519
$ws.Range("T5").Value = 0.6145646077002448

# Row 6 (target cluster: Resolving-Mac) - recomputed TPM values
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03348933333333334
$ws.Range("N6").Value = 0.100468
$ws.Range("O6").Value = 0.003999268915745247
$ws.Range("P6").Value = 0.005404008943283026
$ws.Range("Q6").Value = 0.004060029092666667
$ws.Range("R6").Value = 0.024360174556
$ws.Range("S6").Value = 0.003999268915745247
$ws.Range("T6").Value = 0.005404008943283026
